# Weekly fruit/vegetable price update:
# Insert a new data row at row 251 (pushing the existing rows 251-366 down
# to 252-367) and populate the new row with this week's reading.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(251).Insert()

$ws.Cells.Item(251, 1).Value = 4
$ws.Cells.Item(251, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(251, 3).Value = "Los Lagos"
$ws.Cells.Item(251, 4).Value = 44917
$ws.Cells.Item(251, 5).Value = 10
$ws.Cells.Item(251, 6).Value = 100112017
$ws.Cells.Item(251, 7).Value = "Apio"
$ws.Cells.Item(251, 8).Value = "Americana (o)"
$ws.Cells.Item(251, 9).Value = "Primera"
$ws.Cells.Item(251, 10).Value = 25
$ws.Cells.Item(251, 11).Value = 13000
$ws.Cells.Item(251, 12).Value = 13000
$ws.Cells.Item(251, 13).Value = 13000
$ws.Cells.Item(251, 14).Value = "`$/docena de matas"
$ws.Cells.Item(251, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(251, 16).Value = 2167
$ws.Cells.Item(251, 17).Value = 6
$ws.Cells.Item(251, 18).Value = "Hortaliza"
